$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

# Row 10
$ws.Range("A10").Value = 111936893
$ws.Range("B10").Value = 77650
$ws.Range("E10").Value = 6425
$ws.Range("F10").Value = 'Garnlav'
$ws.Range("G10").Value = 'Alectoria sarmentosa'
$ws.Range("H10").Value = '(Ach.) Ach.'
$ws.Range("Q10").Value = 448742
$ws.Range("R10").Value = 7087502

# Row 11
$ws.Range("A11").Value = 111936869
$ws.Range("B11").Value = 89571
$ws.Range("D11").Value = 'NT'
$ws.Range("E11").Value = 5432
$ws.Range("F11").Value = 'Granticka'
$ws.Range("G11").Value = 'Porodaedalea chrysoloma'
$ws.Range("H11").Value = '(Fr.) Fiasson & Niemelä'
$ws.Range("Q11").Value = 449144
$ws.Range("R11").Value = 7087118

# Row 12
$ws.Range("A12").Value = 111936798
$ws.Range("B12").Value = 56430
$ws.Range("E12").Value = 100109
$ws.Range("F12").Value = 'Tretåig hackspett'
$ws.Range("G12").Value = 'Picoides tridactylus'
$ws.Range("H12").Value = '(Linnaeus, 1758)'
$ws.Range("Q12").Value = 448923
$ws.Range("R12").Value = 7087371
$ws.Range("K12").Value = "'"
$ws.Range("K12").Style = "Normal"
$ws.Range("L12").Value = "'"
$ws.Range("L12").Style = "Normal"
$ws.Range("M12").Value = "'"
$ws.Range("M12").Style = "Normal"
$ws.Range("N12").Value = "'"
$ws.Range("N12").Style = "Normal"
$ws.Range("AC12").Value = 'ringhack äldre'

# Row 13
$ws.Range("A13").Value = 111936868
$ws.Range("B13").Value = 89571
$ws.Range("E13").Value = 5432
$ws.Range("F13").Value = 'Granticka'
$ws.Range("G13").Value = 'Porodaedalea chrysoloma'
$ws.Range("H13").Value = '(Fr.) Fiasson & Niemelä'
$ws.Range("Q13").Value = 448988
$ws.Range("R13").Value = 7087187

# Row 15
$ws.Range("A15").Value = 111936866
$ws.Range("B15").Value = 89571
$ws.Range("Q15").Value = 448766
$ws.Range("R15").Value = 7087417

# Row 16
$ws.Range("A16").Value = 111936870
$ws.Range("B16").Value = 89571
$ws.Range("E16").Value = 5432
$ws.Range("F16").Value = 'Granticka'
$ws.Range("G16").Value = 'Porodaedalea chrysoloma'
$ws.Range("H16").Value = '(Fr.) Fiasson & Niemelä'
$ws.Range("Q16").Value = 449019
$ws.Range("R16").Value = 7087277
$ws.Range("K16").ClearContents()
$ws.Range("L16").ClearContents()
$ws.Range("M16").ClearContents()
$ws.Range("N16").ClearContents()
$ws.Range("AC16").ClearContents()

# Row 17
$ws.Range("A17").Value = 111936865
$ws.Range("B17").Value = 89571
$ws.Range("E17").Value = 5432
$ws.Range("F17").Value = 'Granticka'
$ws.Range("G17").Value = 'Porodaedalea chrysoloma'
$ws.Range("H17").Value = '(Fr.) Fiasson & Niemelä'
$ws.Range("Q17").Value = 448738
$ws.Range("R17").Value = 7087426
$ws.Range("K17").ClearContents()
$ws.Range("L17").ClearContents()
$ws.Range("M17").ClearContents()
$ws.Range("N17").ClearContents()
$ws.Range("AC17").ClearContents()

# Row 18
$ws.Range("B18").Value = 89993

# Row 19
$ws.Range("A19").Value = 111936795
$ws.Range("B19").Value = 56430
$ws.Range("E19").Value = 100109
$ws.Range("F19").Value = 'Tretåig hackspett'
$ws.Range("G19").Value = 'Picoides tridactylus'
$ws.Range("H19").Value = '(Linnaeus, 1758)'
$ws.Range("Q19").Value = 448749
$ws.Range("R19").Value = 7087422
$ws.Range("K19").Value = "'"
$ws.Range("K19").Style = "Normal"
$ws.Range("L19").Value = "'"
$ws.Range("L19").Style = "Normal"
$ws.Range("M19").Value = "'"
$ws.Range("M19").Style = "Normal"
$ws.Range("N19").Value = "'"
$ws.Range("N19").Style = "Normal"
$ws.Range("AC19").Value = 'ringhack äldre'

# Row 20
$ws.Range("A20").Value = 111936792
$ws.Range("B20").Value = 90235
$ws.Range("D20").Value = 'LC'
$ws.Range("E20").Value = 3298
$ws.Range("F20").Value = 'Trådticka'
$ws.Range("G20").Value = 'Climacocystis borealis'
$ws.Range("H20").Value = '(Fr.) Kotl. & Pouzar'
$ws.Range("Q20").Value = 448761
$ws.Range("R20").Value = 7087579

# Row 21
$ws.Range("A21").Value = 111936867
$ws.Range("B21").Value = 89571
$ws.Range("Q21").Value = 448792
$ws.Range("R21").Value = 7087386
